# Auto-generated edit script applying the cell-value changes described in the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.891.39"
$ws.Range("E2").Value = "  -0.33%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.116.66"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "347.67"
$ws.Range("E5").Value = "  -0.15%  "
$ws.Range("E6").Value = "  +0.50%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5195"
$ws.Range("E7").Value = "  +0.76%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.4448"
$ws.Range("E8").Value = "  +0.53%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "54.23"
$ws.Range("E9").Value = "  +3.24%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.09328"
$ws.Range("E10").Value = "  +0.71%  "
$ws.Range("E11").Value = "  +0.83%  "
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "8.406"
$ws.Range("E13").Value = "  +2.31%  "
$ws.Range("B14").Value = "WrappedEther"
$ws.Range("C14").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "2.094.22"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("B15").Value = "Polkadot"
$ws.Range("C15").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "6.857"
$ws.Range("E15").Value = "  +1.82%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "102.30"
$ws.Range("E16").Value = "  +3.43%  "
$ws.Range("E17").Value = "  +0.84%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.008"
$ws.Range("E18").Value = "  +0.48%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "21.45"
$ws.Range("E19").Value = "  +3.69%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.06677"
$ws.Range("E20").Value = "  +0.21%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.294"
$ws.Range("E21").Value = "  +1.04%  "
$ws.Range("E22").Value = "  +0.40%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "29.918.44"
$ws.Range("E23").Value = "  -0.59%  "
$ws.Range("E24").Value = "  +0.47%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.328"
$ws.Range("E25").Value = "  -0.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.351.41"
$ws.Range("E26").Value = "  -0.05%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "22.09"
$ws.Range("E27").Value = "  +0.54%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.560"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "162.31"
$ws.Range("E29").Value = "  +0.09%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "133.98"
$ws.Range("E30").Value = "  +0.45%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.154"
$ws.Range("E31").Value = "  -1.65%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.791"
$ws.Range("E32").Value = "  +8.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.1055"
$ws.Range("E33").Value = "  -0.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "6.242"
$ws.Range("E34").Value = "  +0.21%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "3.976"
$ws.Range("E35").Value = "  +0.48%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.480"
$ws.Range("E36").Value = "  +5.62%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "10.83"
$ws.Range("E37").Value = "  +6.53%  "
$ws.Range("E38").Value = "  +1.07%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.06831"
$ws.Range("E39").Value = "  +0.68%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.7027"
$ws.Range("E40").Value = "  +2.00%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "12.67"
$ws.Range("E41").Value = "  +1.11%  "
$ws.Range("B42").Value = "Algorand"
$ws.Range("C42").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.2253"
$ws.Range("E42").Value = "  -1.10%  "
$ws.Range("B43").Value = "TrustWalletToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.333"
$ws.Range("E43").Value = "  +1.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.6806"
$ws.Range("E44").Value = "  +3.24%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "14.45"
$ws.Range("E45").Value = "  +1.73%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.346"
$ws.Range("E46").Value = "  +2.77%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.006"
$ws.Range("E47").Value = "  +0.50%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.00000000359"
$ws.Range("E48").Value = "  +1.01%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "3.638"
$ws.Range("E49").Value = "  +0.06%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.237"
$ws.Range("E50").Value = "  +6.02%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.222"
$ws.Range("E51").Value = "  +0.19%  "
